$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at 1207 (pushes the existing rows 1207-1269 down to
# 1208-1270, Excel auto-extends the used range / dimension accordingly).
$ws.Rows.Item(1207).Insert()

# Populate the newly inserted row with the new weekly price-point record.
$ws.Range("A1207").Value = 6
$ws.Range("B1207").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1207").Value = "Metropolitana"
$ws.Range("D1207").Value = 45041
$ws.Range("E1207").Value = 13
$ws.Range("F1207").Value = 100112031
$ws.Range("G1207").Value = "Poroto verde"
$ws.Range("H1207").Value = "Magnum"
$ws.Range("I1207").Value = "Primera"
$ws.Range("J1207").Value = 580
$ws.Range("K1207").Value = 20000
$ws.Range("L1207").Value = 22000
$ws.Range("M1207").Value = 20897
$ws.Range("N1207").Value = "$/saco 25 kilos"
$ws.Range("O1207").Value = "Región Metropolitana"
$ws.Range("P1207").Value = 836
$ws.Range("Q1207").Value = 25
$ws.Range("R1207").Value = "Hortaliza"
